$d = $word.ActiveDocument

# Locate the paragraph ending in "LOB1036: Geometria Analítica (Requisito fraco)".
# The deletion starts right after this paragraph's mark (i.e. the paragraph itself
# is left untouched).
$r1 = $d.Content
$found1 = $r1.Find.Execute("LOB1036: Geometria Analítica (Requisito fraco)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPoint = $r1.Paragraphs(1).Range.End

# Locate the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph; the block to be
# removed is: the blank paragraph before it, this paragraph itself, the blank
# paragraph after it, and the page-break paragraph that follows that one. The two
# final (blank / page-break) paragraphs that follow are left in place.
$r2 = $d.Content
$found2 = $r2.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$verParagraph = $r2.Paragraphs(1)
$afterPageBreakParagraph = $verParagraph.Next().Next()
$endPoint = $afterPageBreakParagraph.Range.End

$d.Range($startPoint, $endPoint).Delete()
